$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 667, shifting existing rows 667:712 down to 668:713
$ws.Rows.Item(667).Insert()

# Copy formatting (date style) from the row above into the new row's D cell
$ws.Range("D666").Copy()
$ws.Range("D667").PasteSpecial(-4122) # xlPasteFormats

# Populate the new row 667 with the new record's data
$ws.Cells.Item(667, 1).Value = 8
$ws.Cells.Item(667, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(667, 3).Value = "Coquimbo"
$ws.Cells.Item(667, 4).Value = 45223
$ws.Cells.Item(667, 5).Value = 4
$ws.Cells.Item(667, 6).Value = 100112017
$ws.Cells.Item(667, 7).Value = "Apio"
$ws.Cells.Item(667, 8).Value = "Americana (o)"
$ws.Cells.Item(667, 9).Value = "Primera"
$ws.Cells.Item(667, 10).Value = 1200
$ws.Cells.Item(667, 11).Value = 7000
$ws.Cells.Item(667, 12).Value = 8000
$ws.Cells.Item(667, 13).Value = 7500
$ws.Cells.Item(667, 14).Value = "$/docena de matas"
$ws.Cells.Item(667, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(667, 16).Value = 1250
$ws.Cells.Item(667, 17).Value = 6
$ws.Cells.Item(667, 18).Value = "Hortaliza"
